# Commit: fixed bold/italic/underline, sub and superscript
#
# 1) Mark "excercise" as a flagged spell-check word (spellStart/spellEnd)
#    in the "Why do we need to do so much excercise?" question.
# 2) Turn the "This is another example question?" question (and the 4
#    plain-text answer paragraphs that followed it) into a bold / italic /
#    underline formatting demo, with 3 follow-up demo paragraphs.
# 3) Turn "Another question to pad things out?" into a subscript /
#    superscript formatting demo (with spell-check flags on the made-up
#    words).
$d = $word.ActiveDocument

# --- Change 1 -------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Why do we need to do so much excercise?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$prng = $rng.Paragraphs(1).Range
$xmlFrag0 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="199253C5" w14:textId="0C8D21BC" w:rsidR="00CB3C33" w:rsidRDefault="009C2981"><w:r><w:t>[1 mark]</w:t></w:r><w:r w:rsidR="00CB3C33"><w:t xml:space="preserve">Why do we need to do so much </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00621F9A"><w:t>excercise</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00CB3C33"><w:t>?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$prng.InsertXML($xmlFrag0)

# --- Change 2 -------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("This is another example question?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p1start = $rng.Paragraphs(1).Range.Start
$rng2 = $d.Content
$rng2.Find.Execute("That will come!", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p2end = $rng2.Paragraphs(1).Range.End
$full = $d.Range($p1start, $p2end)
$xmlFrag1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="16BFD43E" w14:textId="35A8F4AC" w:rsidR="00CB3C33" w:rsidRDefault="009C2981"><w:r><w:t xml:space="preserve">[2 marks] </w:t></w:r><w:r><w:t xml:space="preserve">Let’s test </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>bold</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>italics</w:t></w:r><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>underline</w:t></w:r><w:r><w:t>?</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Bold</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Italics</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Underline</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$full.InsertXML($xmlFrag1)

# --- Change 3 -------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Another question to pad things out?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$prng = $rng.Paragraphs(1).Range
$xmlFrag2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="5C1FDD0D" w14:textId="2C4D9C44" w:rsidR="00673B84" w:rsidRDefault="009C2981"><w:r><w:t>[2 marks]</w:t></w:r><w:r w:rsidR="00673B84"><w:t xml:space="preserve">Let’s </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00673B84"><w:t>test</w:t></w:r><w:r w:rsidR="00673B84"><w:rPr><w:vertAlign w:val="subscript"/></w:rPr><w:t>subscript</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00673B84"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00673B84"><w:t>and</w:t></w:r><w:r w:rsidR="00673B84"><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>superscript</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00673B84"><w:t>?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$prng.InsertXML($xmlFrag2)
